$wb = $excel.ActiveWorkbook

# --- OrderItems sheet: selection moves to A2:A15, scrolled down a bit ---
$wsItems = $wb.Worksheets.Item("OrderItems")
$wsItems.Range("A2:A15").Select()

# --- orderIdwithDate sheet becomes the new "homepage": headers switch from the
#     Chinese display names to the English item codes used elsewhere in the
#     workbook (OrderItems!A2:A15), so it can be used to jump/link into the
#     other three pages. ---
$wsDate = $wb.Worksheets.Item("orderIdwithDate")
$wsDate.Range("B1").Value = "quanYaTwoEat"
$wsDate.Range("C1").Value = "quanJiaTwoEatSpicy"
$wsDate.Range("D1").Value = "banYaTwoEat"
$wsDate.Range("E1").Value = "banYaTwoEatSpicy"
$wsDate.Range("F1").Value = "quanYaChopFry"
$wsDate.Range("G1").Value = "quanYaChopFrySpicy"
$wsDate.Range("H1").Value = "banYaChopFry"
$wsDate.Range("I1").Value = "banYaChopFrySpicy"
$wsDate.Range("J1").Value = "quanYaChopPlate"
$wsDate.Range("K1").Value = "banYaChopPlate"
$wsDate.Range("L1").Value = "quanJiShouPaJi"
$wsDate.Range("M1").Value = "banJiShouPaJi"
$wsDate.Range("N1").Value = "heYeBing"
$wsDate.Range("O1").Value = "tianMianJiang"

# This sheet becomes the new active/home tab.
$wsDate.Range("M14").Select()
$wsDate.Activate()
